# Yting_260321 / Ratios.xlsx
# - added support for string labnrs: "Lab. #" column now holds alphanumeric
#   sample labels (e.g. "10815a") instead of bare numeric lab numbers.
# - output of tailing.xlsx / intensities.xlsx + numba speed-up of the age
#   calculation slightly changed the downstream ratio columns (P..U).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New string-typed lab numbers and refreshed ratio values (columns P-U) per row.
$rowUpdates = @(
    @{ Row = 2;  Lab = "10815a"; P = 1.318826771892574;    Q = 0.5335004453826736;   R = 0.1702371670050599;    S = 1.612895259779135;   T = 0.2257985771763353;    U = 1.726512574742643 },
    @{ Row = 3;  Lab = "11069b"; P = 0.03525043890623247;  Q = 5.705288904693922;    R = 0.3294613250209388;    S = 1.213265190703177;   T = 0.01157264535611715;   U = 2.453371522281945 },
    @{ Row = 4;  Lab = "10815a"; P = 1.318759177757272;    Q = 0.6996313538309004;   R = 0.1720259765033036;    S = 1.288866266384345;   T = 0.2271470050020687;    U = 1.004207089032775 },
    @{ Row = 5;  Lab = "11070a"; P = 0.005921825770778144; Q = 2.813032265449945;    R = 0.2259939097674143;    S = 1.381875981856847;   T = 0.001274691673562648;  U = 2.97747271498213 },
    @{ Row = 6;  Lab = "10815a"; P = 1.313374031600917;    Q = 0.5419107723895897;   R = 0.1815114725115821;    S = 1.276965247862709;   T = 0.2393293430170869;    U = 1.28265146190901 },
    @{ Row = 7;  Lab = "11071a"; P = 0.01783305396177997;  Q = 1.001748480524192;    R = 0.249544100283467;     S = 0.5610059207786317;  T = 0.004458367135931112;  U = 0.6951876357866389 },
    @{ Row = 8;  Lab = "10815a"; P = 1.313892579368417;    Q = 0.4004394967171863;   R = 0.1792768137488384;    S = 1.36931034548719;    T = 0.2374191339979786;    U = 1.300845587375957 },
    @{ Row = 9;  Lab = "11072a"; P = 0.1438931916272418;   Q = 0.365815313729743;    R = 0.01854280137926809;   S = 0.2678043872163008;  T = 0.002664781889775179;  U = 0.3003130595442578 },
    @{ Row = 10; Lab = "10815a"; P = 1.312935307459668;    Q = 0.5991393735161512;   R = 0.1781077787529703;    S = 1.347426895188685;   T = 0.2340798179609642;    U = 1.249551913334694 },
    @{ Row = 11; Lab = "11074b"; P = 0.1434984853044898;   Q = 0.3364171721037157;   R = 0.0009631928586580701; S = 0.255338915736973;   T = 0.0001387130736508046; U = 0.2564753741715639 },
    @{ Row = 12; Lab = "10815a"; P = 1.318059551970325;    Q = 0.4537665522382115;   R = 0.1784094620679704;    S = 1.295662586723902;   T = 0.2343682381076229;    U = 1.283193429135786 }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.Lab
    $ws.Range("P$r").Value = $u.P
    $ws.Range("Q$r").Value = $u.Q
    $ws.Range("R$r").Value = $u.R
    $ws.Range("S$r").Value = $u.S
    $ws.Range("T$r").Value = $u.T
    $ws.Range("U$r").Value = $u.U
}

# Column U (21) narrowed slightly to match the new layout.
$ws.Columns.Item(21).ColumnWidth = 18.857142857142858
